# "Criado os arquivos de trabalho"
# On the "Planilha2" worksheet (2nd sheet / sheet2.xml):
#  - clear the stray "x" mark that was in E16
#  - add four new task rows (28-31) listing the work files that were
#    created, each marked with an "x" in the corresponding column,
#    copying the formatting used by the previous row (27)
#  - leave the selection on E16, matching the saved workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 16 no longer has the "x" marker in column E
$ws.Cells.Item(16, 5).Value = ""

# New rows describing the created work files
$novosArquivos = @(
    "index.html",
    "style.css(front)",
    "f1_main.js",
    "f1_class.js"
)

$linhaModelo = 27
for ($i = 0; $i -lt $novosArquivos.Length; $i++) {
    $linha = $linhaModelo + 1 + $i

    # copy formatting/values from the row above, then rename column B
    $ws.Range("B$linhaModelo`:E$linhaModelo").Copy($ws.Range("B$linha`:E$linha"))
    $ws.Cells.Item($linha, 2).Value = $novosArquivos[$i]
}

# Restore the selection saved in the workbook (cell E16 on Planilha2)
$ws.Activate()
$ws.Range("E16").Select()
